$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (GitHub Actions refresh)

$ws.Range("D2").Value = "67.142.41"
$ws.Range("E2").Value = "  +4.29%  "

$ws.Range("D3").Value = "3.257.67"
$ws.Range("E3").Value = "  +2.42%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'577.81"
$ws.Range("E5").Value = "  +2.11%  "

$ws.Range("D6").Value = "'178.60"
$ws.Range("E6").Value = "  +4.35%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "'0.603"
$ws.Range("E8").Value = "  -1.14%  "

$ws.Range("D9").Value = "3.252.72"
$ws.Range("E9").Value = "  +2.39%  "

$ws.Range("E10").Value = "  +3.80%  "

$ws.Range("D11").Value = "'6.76"
$ws.Range("E11").Value = "  +1.70%  "

$ws.Range("D12").Value = "'0.414"
$ws.Range("E12").Value = "  +4.17%  "

$ws.Range("D13").Value = "3.817.24"
$ws.Range("E13").Value = "  +2.35%  "

$ws.Range("E14").Value = "  +0.76%  "

$ws.Range("D15").Value = "'28.23"
$ws.Range("E15").Value = "  +2.85%  "

$ws.Range("D16").Value = "67.092.62"
$ws.Range("E16").Value = "  +4.24%  "

$ws.Range("E17").Value = "  +2.57%  "

$ws.Range("D18").Value = "3.253.37"
$ws.Range("E18").Value = "  +2.11%  "

$ws.Range("D19").Value = "'5.87"
$ws.Range("E19").Value = "  +2.00%  "

$ws.Range("E20").Value = "  +3.08%  "

$ws.Range("D21").Value = "'375.14"
$ws.Range("E21").Value = "  +6.01%  "

$ws.Range("E22").Value = "  +6.22%  "

$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("D24").Value = "'70.96"
$ws.Range("E24").Value = "  +2.71%  "

$ws.Range("D25").Value = "'0.513"
$ws.Range("E25").Value = "  +1.75%  "

$ws.Range("D26").Value = "3.395.70"
$ws.Range("E26").Value = "  +2.67%  "

$ws.Range("E27").Value = "  -1.45%  "

$ws.Range("D28").Value = "'9.89"
$ws.Range("E28").Value = "  +3.30%  "

$ws.Range("E29").Value = "  +1.76%  "

$ws.Range("E30").Value = "  +0.74%  "

$ws.Range("E31").Value = "  +3.51%  "

$ws.Range("D32").Value = "'5.65"
$ws.Range("E32").Value = "  +0.41%  "

$ws.Range("D33").Value = "'22.62"

$ws.Range("E34").Value = "  +0.03%  "

$ws.Range("D35").Value = "'1.27"
$ws.Range("E35").Value = "  +4.57%  "

$ws.Range("E36").Value = "  +2.55%  "

$ws.Range("D37").Value = "'166.77"
$ws.Range("E37").Value = "  +7.39%  "

$ws.Range("D38").Value = "'1.51"
$ws.Range("E38").Value = "  +4.26%  "

$ws.Range("D39").Value = "'0.857"
$ws.Range("E39").Value = "  +4.91%  "

$ws.Range("E40").Value = "  +10.03%  "

$ws.Range("D41").Value = "'27.09"
$ws.Range("E41").Value = "  +4.49%  "

$ws.Range("D42").Value = "'2.59"
$ws.Range("E42").Value = "  +1.00%  "

$ws.Range("D43").Value = "2.770.30"
$ws.Range("E43").Value = "  +5.56%  "

$ws.Range("D44").Value = "'6.47"
$ws.Range("E44").Value = "  +7.58%  "

$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").Value = "'355.58"
$ws.Range("E45").Value = "  +8.93%  "

$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").Value = "'4.40"
$ws.Range("E46").Value = "  +4.72%  "

$ws.Range("D47").Value = "'25.49"
$ws.Range("E47").Value = "  +5.93%  "

$ws.Range("E48").Value = "  +1.90%  "

$ws.Range("E49").Value = "  +2.27%  "

$ws.Range("D50").Value = "'0.0281"
$ws.Range("E50").Value = "  +3.50%  "

$ws.Range("E51").Value = "  +0.80%  "

